$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.054.84"
$ws.Range("E2").Value = "  +4.07%  "

$ws.Range("D3").Value = "2.463.87"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.53"
$ws.Range("E5").Value = "  +3.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.23"
$ws.Range("E6").Value = "  +10.23%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.512"
$ws.Range("E8").Value = "  +2.00%  "

$ws.Range("D9").Value = "2.460.23"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  +5.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  +3.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  +3.69%  "

$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").Value = "2.886.54"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "57.262.75"
$ws.Range("E15").Value = "  +4.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.05"
$ws.Range("E16").Value = "  +3.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +3.60%  "

$ws.Range("D18").Value = "2.465.34"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  +5.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "326.69"
$ws.Range("E20").Value = "  +4.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  +3.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.95"
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("D28").Value = "2.560.19"
$ws.Range("E28").Value = "  +0.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.55"
$ws.Range("E29").Value = "  +3.95%  "

$ws.Range("D30").Value = "0.0₃0819"
$ws.Range("E30").Value = "  +6.72%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.45"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.53"
$ws.Range("E33").Value = "  +4.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.20"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +2.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("E36").Value = "  +7.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("E37").Value = "  +2.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.78"
$ws.Range("E38").Value = "  +5.99%  "

$ws.Range("E39").Value = "  +9.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.26"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.51"
$ws.Range("E41").Value = "  +2.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0559"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.607"
$ws.Range("E44").Value = "  +1.97%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.95"
$ws.Range("E45").Value = "  +5.57%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0960"
$ws.Range("E46").Value = "  +8.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "266.85"
$ws.Range("E47").Value = "  +4.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +3.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.25"
$ws.Range("E49").Value = "  +1.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.83"
$ws.Range("E50").Value = "  +4.58%  "

$ws.Range("E51").Value = "  +27.31%  "
